$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 93, shifting existing rows
# 93-102 down to 94-103 (their values are untouched by this operation).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly price record.
$ws.Cells.Item(93, 1).Value2()  = 8
$ws.Cells.Item(93, 2).Value2()  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(93, 3).Value2()  = "Coquimbo"
$ws.Cells.Item(93, 4).Value2()  = 44918
$ws.Cells.Item(93, 5).Value2()  = 4
$ws.Cells.Item(93, 6).Value2()  = 100112030
$ws.Cells.Item(93, 7).Value2()  = "Poroto granado"
$ws.Cells.Item(93, 8).Value2()  = "Sin especificar"
$ws.Cells.Item(93, 9).Value2()  = "Primera"
$ws.Cells.Item(93, 10).Value2() = 470
$ws.Cells.Item(93, 11).Value2() = 34000
$ws.Cells.Item(93, 12).Value2() = 35000
$ws.Cells.Item(93, 13).Value2() = 34500
$ws.Cells.Item(93, 14).Value2() = "$/malla 25 kilos"
$ws.Cells.Item(93, 15).Value2() = "Provincia de Limarí"
$ws.Cells.Item(93, 16).Value2() = 1380
$ws.Cells.Item(93, 17).Value2() = 25
$ws.Cells.Item(93, 18).Value2() = "Hortaliza"

# Make sure the date cell keeps the workbook's date number format (style
# index 2) used by the rest of column D.
$ws.Cells.Item(93, 4).NumberFormat = $ws.Cells.Item(94, 4).NumberFormat
